$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update optimization results for less-constrained firms (rows 52-60)
# Row 52
$ws.Range("N52").Value = 0.011675815967
$ws.Range("O52").Value = 0.010988030018
$ws.Range("P52").Value = 1.693066138998
$ws.Range("Q52").Value = 1.5309694499449999
$ws.Range("R52").Value = 120.811907979879
$ws.Range("S52").Value = 4.2763646614770003
$ws.Range("T52").Value = 5.0515865345339996
$ws.Range("U52").Value = 2.229535800286
$ws.Range("V52").Value = 0.38668183671599998
$ws.Range("W52").Value = 0.63268452493899996
$ws.Range("X52").Value = 0.81665925758699998

# Row 53
$ws.Range("N53").Value = 0.043941017777999997
$ws.Range("O53").Value = 0.041054243781999997
$ws.Range("P53").Value = 1.546434727062
$ws.Range("Q53").Value = 1.399208154439
$ws.Range("R53").Value = 12.528886566156
$ws.Range("S53").Value = 1.0223204031390001
$ws.Range("T53").Value = 1.264416426413
$ws.Range("U53").Value = 4.0942748305540002
$ws.Range("V53").Value = 0.30398230445899999
$ws.Range("W53").Value = 0.442374368465
$ws.Range("X53").Value = 0.56247309236900001

# Row 54
$ws.Range("N54").Value = 0.055791661999999999
$ws.Range("O54").Value = 0.052010193203000002
$ws.Range("P54").Value = 1.513096033591
$ws.Range("Q54").Value = 1.3737313002260001
$ws.Range("R54").Value = 4.4937827015109999
$ws.Range("S54").Value = 0.55174393691300005
$ws.Range("T54").Value = 0.68463885550299997
$ws.Range("U54").Value = 5.7605952186779996
$ws.Range("V54").Value = 0.27309294396900002
$ws.Range("W54").Value = 0.351765876496
$ws.Range("X54").Value = 0.44372465680399997

# Row 55
$ws.Range("N55").Value = 0.0081897672119999997
$ws.Range("O55").Value = 0.00791205377
$ws.Range("P55").Value = 1.7514860972270001
$ws.Range("Q55").Value = 1.585021037338
$ws.Range("R55").Value = 644.29736053513295
$ws.Range("S55").Value = 20.647616166904001
$ws.Range("T55").Value = 23.825868158281001
$ws.Range("U55").Value = 1.906389575583
$ws.Range("V55").Value = 0.43656936243700001
$ws.Range("W55").Value = 0.65482685345299996
$ws.Range("X55").Value = 0.84637348485899999

# Row 56
$ws.Range("N56").Value = 0.040212931409000001
$ws.Range("O56").Value = 0.038497365759999999
$ws.Range("P56").Value = 1.5768812956239999
$ws.Range("Q56").Value = 1.425203107018
$ws.Range("R56").Value = 38.828523486845
$ws.Range("S56").Value = 2.6391529445950002
$ws.Range("T56").Value = 3.2080107949590002
$ws.Range("U56").Value = 3.2895803722610002
$ws.Range("V56").Value = 0.34188854248200001
$ws.Range("W56").Value = 0.48888964209000002
$ws.Range("X56").Value = 0.62304656269600001

# Row 57
$ws.Range("N57").Value = 0.053802480419000001
$ws.Range("O57").Value = 0.051333357088999998
$ws.Range("P57").Value = 1.530615660117
$ws.Range("Q57").Value = 1.388551751182
$ws.Range("R57").Value = 11.160036264612
$ws.Range("S57").Value = 1.1157244378169999
$ws.Range("T57").Value = 1.374698045775
$ws.Range("U57").Value = 4.4807277071409999
$ws.Range("V57").Value = 0.30646791229300002
$ws.Range("W57").Value = 0.40227830299099998
$ws.Range("X57").Value = 0.50886283933599996

# Row 58
$ws.Range("N58").Value = 0.0050722943560000001
$ws.Range("O58").Value = 0.0050848720279999997
$ws.Range("P58").Value = 1.8291973798329999
$ws.Range("Q58").Value = 1.6625538060559999
$ws.Range("R58").Value = 6646.7809369616198
$ws.Range("S58").Value = 195.953452813708
$ws.Range("T58").Value = 219.41810186734401
$ws.Range("U58").Value = 1.664115735917
$ws.Range("V58").Value = 0.48827158859300002
$ws.Range("W58").Value = 0.67023641358499997
$ws.Range("X58").Value = 0.86726488241800004

# Row 59
$ws.Range("N59").Value = 0.034714399745999998
$ws.Range("O59").Value = 0.034400825147999997
$ws.Range("P59").Value = 1.615903268594
$ws.Range("Q59").Value = 1.4593310582929999
$ws.Range("R59").Value = 178.110285268025
$ws.Range("S59").Value = 10.180519812041
$ws.Range("T59").Value = 12.184290868039
$ws.Range("U59").Value = 2.6903158978059998
$ws.Range("V59").Value = 0.38428085230800002
$ws.Range("W59").Value = 0.53094029280300004
$ws.Range("X59").Value = 0.67849427514399996

# Row 60
$ws.Range("N60").Value = 0.049427278024000001
$ws.Range("O60").Value = 0.048730297972
$ws.Range("P60").Value = 1.5623064079919999
$ws.Range("Q60").Value = 1.4148099488749999
$ws.Range("R60").Value = 37.236168857734
$ws.Range("S60").Value = 3.0142100131610001
$ws.Range("T60").Value = 3.6133512209499998
$ws.Range("U60").Value = 3.5656139495630002
$ws.Range("V60").Value = 0.34104787539300002
$ws.Range("W60").Value = 0.45378575776699998
$ws.Range("X60").Value = 0.57524149968899996

# Update the view state: frozen pane top-left + active selection
$ws.Range("A60").Select()

Write-Output "done"
